$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '26.494.45'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.848.78'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '262.60'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -6.89%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5135'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.25%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3213'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -9.22%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06769'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.40%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.04'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.60%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.7674'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -6.13%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07691'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").Value = '1.854.83'
$ws.Range("E13").Value = '  -0.79%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '88.79'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.034'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.96%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  -2.21%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007893'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").Value = '26.522.06'
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '2.091.45'
$ws.Range("E21").Value = '  -0.57%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.559'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -5.19%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.540'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -5.78%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.962'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -4.46%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.336'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '144.90'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.654'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '16.98'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.22%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '111.46'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.55%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.213'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -4.26%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.169'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.31%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.08734'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.87%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04841'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.138'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.36%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.848'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.49%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6912'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -7.53%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.115'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -5.22%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01804'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.95%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.213'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -8.58%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.4920'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -6.81%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '113.64'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.75%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.9043'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -6.96%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '6.164'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.13%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '7.791'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.75%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4251'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -7.45%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.1267'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -7.05%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.164'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.05897'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '35.02'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.03%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.423'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -6.13%  '
